$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '41.665.42'
$ws.Range("E2").Value = '  -0.19%  '

$ws.Range("D3").Value = '2.472.90'
$ws.Range("E3").Value = '  +0.20%  '

$ws.Range("E4").Value = '  +0.15%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '320.10'
$ws.Range("E5").Value = '  +1.29%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '91.84'
$ws.Range("E6").Value = '  -1.17%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.549'
$ws.Range("E7").Value = '  -0.05%  '

$ws.Range("E8").Value = '  +0.08%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.510'
$ws.Range("E9").Value = '  -0.95%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '32.87'
$ws.Range("E10").Value = '  +0.52%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0852'
$ws.Range("E11").Value = '  +1.25%  '

$ws.Range("E12").Value = '  -1.00%  '

$ws.Range("D13").Value = '2.857.75'
$ws.Range("E13").Value = '  +0.36%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '6.87'
$ws.Range("E14").Value = '  -0.09%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '15.44'
$ws.Range("E15").Value = '  -2.05%  '

$ws.Range("D16").Value = '2.476.22'
$ws.Range("E16").Value = '  +1.33%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.790'
$ws.Range("E17").Value = '  +1.28%  '

$ws.Range("D18").Value = '41.656.65'
$ws.Range("E18").Value = '  -0.14%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '6.42'
$ws.Range("E19").Value = '  -1.06%  '

$ws.Range("D20").Value = '0.0₃0939'
$ws.Range("E20").Value = '  -0.77%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '71.35'

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '11.21'
$ws.Range("E22").Value = '  -1.85%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '238.89'
$ws.Range("E23").Value = '  +0.00%  '

$ws.Range("E24").Value = '  +0.97%  '

$ws.Range("E25").Value = '  +0.68%  '

$ws.Range("E26").Value = '  -0.04%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '24.91'
$ws.Range("E27").Value = '  +1.22%  '

$ws.Range("E28").Value = '  -1.27%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '9.73'
$ws.Range("E29").Value = '  -0.55%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '36.53'
$ws.Range("E30").Value = '  +2.75%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '156.94'
$ws.Range("E31").Value = '  +0.77%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '5.41'
$ws.Range("E32").Value = '  -2.08%  '

$ws.Range("E33").Value = '  +0.03%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.0768'
$ws.Range("E34").Value = '  +0.81%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '2.57'
$ws.Range("E35").Value = '  -0.09%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '17.04'
$ws.Range("E36").Value = '  -2.55%  '

$ws.Range("B37").Value = 'Stellar'
$ws.Range("C37").Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.116'
$ws.Range("E37").Value = '  +1.07%  '

$ws.Range("B38").Value = 'ARBITRUM'
$ws.Range("C38").Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '1.83'
$ws.Range("E38").Value = '  +1.92%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '2.86'
$ws.Range("E39").Value = '  -1.06%  '

$ws.Range("E40").Value = '  +0.18%  '

$ws.Range("E41").Value = '  -0.08%  '

$ws.Range("E42").Value = '  -2.23%  '

$ws.Range("D43").Value = '1.998.90'

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.0282'
$ws.Range("E44").Value = '  -0.24%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '18.67'
$ws.Range("E45").Value = '  -1.34%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '2.96'
$ws.Range("E46").Value = '  +0.24%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '9.43'
$ws.Range("E47").Value = '  +4.35%  '

$ws.Range("D48").Value = '2.738.06'
$ws.Range("E48").Value = '  +1.22%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '97.44'
$ws.Range("E49").Value = '  +0.50%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '75.90'
$ws.Range("E50").Value = '  +4.84%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '67.08'
$ws.Range("E51").Value = '  -0.03%  '
